# Generate Report for Handoff
# Updates the "In Translation" status to "Ready for handoff" across the
# Overview / zh-cn / de-de sheets, refreshes the associated handoff
# timestamps, and widens the Status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newStatusWidth = 16.333333333333336   # widest achievable approximation of the authored 17.2159881591797 char width

# ---- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus                 # zh-cn status
$wsOverview.Range("F2").Value = $newStatus                 # de-de status
$wsOverview.Range("G2").Value = "2016-08-27 14:39:09"      # Latest HO Xliff Generate Date
$wsOverview.Range("E1").ColumnWidth = $newStatusWidth
$wsOverview.Range("F1").ColumnWidth = $newStatusWidth

# ---- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus                     # Status
$wsZhCn.Range("H2").Value = "2016-08-27 14:39:04"          # Latest Handoff Datetime
$wsZhCn.Range("C1").ColumnWidth = $newStatusWidth

# ---- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus                     # Status
$wsDeDe.Range("H2").Value = "2016-08-27 14:39:09"          # Latest Handoff Datetime
$wsDeDe.Range("C1").ColumnWidth = $newStatusWidth
